$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Non-breaking space used after "Eoin Morgan" in the existing data
$nbsp = [char]0x00A0
$batsman = "Eoin Morgan" + $nbsp

# Insert two new columns (ownTeam, oppTeam) before the existing "batsman"
# column (old D), shifting old D:I -> F:K.
$ws.Range("D1:E1").EntireColumn.Insert()

# --- Header row ---
$ws.Range("D1").Value = "ownTeam"
$ws.Range("E1").Value = "oppTeam"

# --- Existing rows: add ownTeam / oppTeam values ---
$ws.Range("D2").Value = "Kolkata Knight Riders"
$ws.Range("E2").Value = "Chennai Super Kings"

$ws.Range("D3").Value = "Kolkata Knight Riders"
$ws.Range("E3").Value = "Rajasthan Royals"

$ws.Range("D4").Value = "Kolkata Knight Riders"
$ws.Range("E4").Value = "Kings XI Punjab"

# --- New rows 5-8 ---
# All values are written as text (matching the rest of the sheet, where
# even numeric-looking values are stored as strings), so force the
# number format to Text before assigning to avoid Excel's automatic
# number conversion.
$newRows = @(
  @{ row=5; A=" Abu Dhabi"; B=" September 26 2020"; C="KKR won by 7 wickets (with 12 balls remaining)"; D="Kolkata Knight Riders"; E="Sunrisers Hyderabad"; G="42"; H="29"; I="3"; J="2"; K="144.82" },
  @{ row=6; A=" Sharjah"; B=" October 03 2020"; C="Capitals won by 18 runs"; D="Kolkata Knight Riders"; E="Delhi Capitals"; G="44"; H="18"; I="1"; J="5"; K="244.44" },
  @{ row=7; A=" Abu Dhabi"; B=" September 23 2020"; C="Mumbai won by 49 runs"; D="Kolkata Knight Riders"; E="Mumbai Indians"; G="16"; H="20"; I="1"; J="1"; K="80.00" },
  @{ row=8; A=" Sharjah"; B=" October 12 2020"; C="RCB won by 82 runs"; D="Kolkata Knight Riders"; E="Royal Challengers Bangalore"; G="8"; H="12"; I="1"; J="0"; K="66.66" }
)

# Set the number format to Text up-front for the whole numeric block
# (G5:K8) so every one of those cells shares a single style entry
# instead of one per assignment.
$ws.Range("G5:K8").NumberFormat = "@"

foreach ($r in $newRows) {
    $row = $r.row
    $ws.Range("A$row").Value = $r.A
    $ws.Range("B$row").Value = $r.B
    $ws.Range("C$row").Value = $r.C
    $ws.Range("D$row").Value = $r.D
    $ws.Range("E$row").Value = $r.E
    $ws.Range("F$row").Value = $batsman

    $ws.Range("G$row").Value = $r.G
    $ws.Range("H$row").Value = $r.H
    $ws.Range("I$row").Value = $r.I
    $ws.Range("J$row").Value = $r.J
    $ws.Range("K$row").Value = $r.K
}
